# Kontroller.xlsx — "Statistics features was added to script"
#
# Applies the control-log updates on the "Controls" sheet:
#   - row 2 (control #1): due date moves 2022-06-09 -> 2022-06-13, marked verified (D2 = "X")
#   - row 3 (control #2): due date moves 2022-05-29 -> 2022-06-29
#   - row 5 (control #4, was "Teleworking"): repurposed into another
#     "Verify Screening processes" pass, due 2022-06-12, marked verified (D5 = "X"),
#     and its date column is restyled to the same format used by rows 2-4 (short date)
#   - two brand new rows (6 & 7, controls #5 & #6) are appended, each another
#     "Verify Screening processes" entry due 2022-06-29, owned by "Knud", styled with
#     the custom DD-MM-YYYY date format that row 5 used to carry
#   - the sheet's selection / active cell ends up on I6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controls")

# --- Propagate the *original* C5 number format (the custom "DD-MM-YYYY" format) onto
# the two brand-new rows before we touch C5's own formatting below.
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C7").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2: push the due date out a few days and mark it verified
$ws.Range("C2").Value = 44725
$ws.Range("D2").Value = "X"

# --- Row 3: push the due date out a month
$ws.Range("C3").Value = 44741

# --- Row 5: turn the old "Teleworking" entry into another screening check, restyle its
# date cell to match the other rows (numeric short-date format, same as C2/C3/C4)
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B5").Value = "Verify Screening processes"
$ws.Range("C5").Value = 44724
$ws.Range("D5").Value = "X"

# --- New row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Verify Screening processes"
$ws.Range("C6").Value = 44741
$ws.Range("E6").Value = "Knud"

# --- New row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Verify Screening processes"
$ws.Range("C7").Value = 44741
$ws.Range("E7").Value = "Knud"

# --- Move the selection/active cell to I6, matching the saved view state
$ws.Range("I6").Select() | Out-Null
